$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 14

# Row 3
$ws.Range("C3").Value = 14
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 14
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 0.9333333333333333
$ws.Range("H4").Value = 1

# Row 5
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 0.875
$ws.Range("H5").Value = 1

# Row 6
$ws.Range("C6").Value = 12
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 14
$ws.Range("H6").Value = 0.8571428571428571

# Row 7
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 0.9333333333333333
$ws.Range("H7").Value = 1

# Row 8
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 14
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = 0

# Row 9
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 12
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 0.5384615384615384
$ws.Range("H9").Value = 1

# Row 10
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 21
$ws.Range("G10").Value = 0.6666666666666666
$ws.Range("H10").Value = 1

# Row 11
$ws.Range("C11").Value = 14
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 16
$ws.Range("G11").Value = 0.875
$ws.Range("H11").Value = 1

# Row 12
$ws.Range("C12").Value = 14
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 14
$ws.Range("H12").Value = 1

# Row 13
$ws.Range("C13").Value = 14
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 14
$ws.Range("H13").Value = 1

# Row 14
$ws.Range("C14").Value = 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 0.9333333333333333
$ws.Range("H14").Value = 1

# Row 15
$ws.Range("C15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 14
$ws.Range("H15").Value = 1

# Row 16
$ws.Range("B16").Value = "invalid_stamp"
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 0.8571428571428571
$ws.Range("H16").Value = 1

# Row 17
$ws.Range("B17").Value = "key"
$ws.Range("C17").Value = 14
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = 1

# Row 18
$ws.Range("B18").Value = "ladder"
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 0.9333333333333333
$ws.Range("H18").Value = 1

# Row 19
$ws.Range("B19").Value = "lock"
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 0.8235294117647058
$ws.Range("H19").Value = 1

# Row 20
$ws.Range("B20").Value = "lotus"
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 0.56
$ws.Range("H20").Value = 1

# Row 21
$ws.Range("B21").Value = "loud_speaker"
$ws.Range("C21").Value = 14
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 14
$ws.Range("H21").Value = 1

# Row 22
$ws.Range("B22").Value = "mother_and_child"
$ws.Range("C22").Value = 9
$ws.Range("F22").Value = 14
$ws.Range("H22").Value = 0.6428571428571429

# Row 23
$ws.Range("B23").Value = "namaste"
$ws.Range("C23").Value = 14
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 14
$ws.Range("H23").Value = 1

# Row 24
$ws.Range("B24").Value = "nepali_big_basket"
$ws.Range("C24").Value = 0
$ws.Range("E24").Value = 14
$ws.Range("F24").Value = 14
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = 0

# Row 25
$ws.Range("B25").Value = "nepali_cap"
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 0.5384615384615384
$ws.Range("H25").Value = 1

# Row 26
$ws.Range("B26").Value = "nepali_jug"
$ws.Range("C26").Value = 14
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 14
$ws.Range("H26").Value = 1

# Row 27
$ws.Range("B27").Value = "nepali_madal"
$ws.Range("C27").Value = 12
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 14
$ws.Range("H27").Value = 0.8571428571428571

# Row 28
$ws.Range("B28").Value = "nepali_small_basket"
$ws.Range("C28").Value = 14
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 14
$ws.Range("H28").Value = 1

# Row 29
$ws.Range("B29").Value = "owl"
$ws.Range("C29").Value = 1
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 14
$ws.Range("H29").Value = 0.07142857142857142

# Row 30
$ws.Range("B30").Value = "pen"
$ws.Range("C30").Value = 11
$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 14
$ws.Range("H30").Value = 0.7857142857142857

# Row 31
$ws.Range("B31").Value = "roaster"
$ws.Range("C31").Value = 14
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 14
$ws.Range("H31").Value = 1

# Row 32
$ws.Range("B32").Value = "sheep"
$ws.Range("C32").Value = 14
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 14
$ws.Range("H32").Value = 1

# Row 33
$ws.Range("C33").Value = 14
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 18
$ws.Range("G33").Value = 0.7777777777777778
$ws.Range("H33").Value = 1

# Row 34
$ws.Range("C34").Value = 14
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 14
$ws.Range("H34").Value = 1

# Row 35
$ws.Range("C35").Value = 14
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 14
$ws.Range("H35").Value = 1

# Row 36
$ws.Range("C36").Value = 13
$ws.Range("D36").Value = 3
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 16
$ws.Range("G36").Value = 0.8125
$ws.Range("H36").Value = 1

# Row 37
$ws.Range("C37").Value = 14
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 15
$ws.Range("G37").Value = 0.9333333333333333
$ws.Range("H37").Value = 1

# Row 38
$ws.Range("C38").Value = 14
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 15
$ws.Range("G38").Value = 0.9333333333333333
$ws.Range("H38").Value = 1

# Row 39
$ws.Range("C39").Value = 13
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 14
$ws.Range("H39").Value = 0.9285714285714286

# Row 40
$ws.Range("C40").Value = 14
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 14
$ws.Range("H40").Value = 1

# Row 41
$ws.Range("B41").Value = "valid_stamp"
$ws.Range("C41").Value = 8
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 8
$ws.Range("H41").Value = 1

# Row 42
$ws.Range("B42").Value = "water_glass"
$ws.Range("C42").Value = 14
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 14
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 1

# Row 43
$ws.Range("B43").Value = "water_jug"
$ws.Range("C43").Value = 14
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 15
$ws.Range("G43").Value = 0.9333333333333333
$ws.Range("H43").Value = 1

# Row 44
$ws.Range("B44").Value = "woman_man"
$ws.Range("C44").Value = 13
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = 14
$ws.Range("H44").Value = 0.9285714285714286

# Row 45 (new row) - copy formatting from row 44, then set values
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "wooden_wheel"
$ws.Range("C45").Value = 14
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 14
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 1
